$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Replace the e-mail addresses stored in column A with their lower-cased
# versions ("IvánGarcía@..." -> "ivangarcía@...", "AllissonFlores@..." ->
# "allissonflores@...", "Eunice@..." -> "eunice@...").
$map = @{
    "IvánGarcía@beeckerco.com"     = "ivangarcía@beeckerco.com"
    "AllissonFlores@beeckerco.com" = "allissonflores@beeckerco.com"
    "Eunice@beeckerco.com"         = "eunice@beeckerco.com"
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}

# Widen the workbook window (matches windowWidth/windowHeight change in
# the saved workbook view; values are twips in the XML, points in the COM
# object model).
$win = $wb.Windows.Item(1)
$win.Width = 29040 / 20
$win.Height = 15840 / 20
